$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
Write-Host "scrollrow before: $($excel.ActiveWindow.ScrollRow)"
$excel.ActiveWindow.ScrollRow = 9
Write-Host "scrollrow after: $($excel.ActiveWindow.ScrollRow)"
$ws.Range("A25").Select()
Write-Host "done"
